$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 ("Objetivos:") gets a real objectives paragraph in B/C -------
$objetivos = "Desenvolver um projeto sobre tema de Engenharia de Produção, similar a situações que os alunos irão encontrar na vida real no efetivo exercício de sua profissão, `nAplicar e integrar conhecimentos adquiridos em demais disciplinas de seu curso`nDesenvolver competências técnicas, as relacionadas ao projeto em si, bem como competências transversais (habilidades e atitudes), num ambiente de aprendizagem baseado em PBL (Project-Baed Learning)."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# --- Insert a brand-new row at 13 ("Docentes responsáveis:" gets its B/C)
# This pushes the old rows 13..21 down to 14..22, which already carries
# their existing content/styles/heights to the right place.
$ws.Rows.Item(13).Insert()

# The freshly inserted row 13 copied formatting from the row above onto
# column A; the target layout has no A13 cell at all, so drop it.
$ws.Range("A13").Clear()

$docente = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente

# --- Row 14 ("Programa resumido:") ---------------------------------------
$resumo = "Tópicos que abordem o tema do projeto de seu planejamento a execução."
$ws.Range("B14").Value = $resumo
$ws.Range("C14").Value = $resumo

# --- Row 16 ("Programa:") -------------------------------------------------
$programa = "Noções de Gestão de Projetos`nOrganização do tempo: dimensão pessoal;`nTécnicas para a realização de apresentações;`nNoções de Aprendizagem Baseada em Projetos`nTrabalho em Grupo, Equipes e times. `nPostura e Ética Profissional`nTécnicas para redação de relatório técnico;`nTutoria de projetos.`nAssuntos Técnicos específicos relacionados com o tema do projeto."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# --- Row 19 ("Método:") ---------------------------------------------------
$metodo = "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.`n`nOs alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. `nCada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.`nAs aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# --- Row 20 ("Critério:") --------------------------------------------------
$criterio = "A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.`nO detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# --- Row 21 ("Norma de recuperação:") --------------------------------------
$recuperacao = "Não há recuperação"
$ws.Range("B21").Value = $recuperacao
$ws.Range("C21").Value = $recuperacao

# --- Row 22 ("Bibliografia:") ------------------------------------------------
$bibliografia = "Artigos sobre metodologias ativas de aprendizagem e  Project Based Learning.`nLivros e Artigos científicos relacionados com o tema do projeto."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
